$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141; this shifts the existing rows 141-145 down to 142-146
$ws.Rows.Item(141).Insert()

# Copy the "static" (unchanged) column values from the row that is now 142
# (these values are identical to what the new row 141 needs) and set the
# style for column D (date format) to match.
$ws.Range("A141").Value = 11
$ws.Range("B141").Value = "Vega Monumental Concepción"
$ws.Range("C141").Value = "Bíobío"
$ws.Range("D141").Style = "Normal"
$ws.Range("D141").NumberFormat = $ws.Range("D142").NumberFormat
$ws.Range("D141").Value = 45223
$ws.Range("E141").Value = 8
$ws.Range("F141").Value = 100112037
$ws.Range("G141").Value = "Cebollín"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 35
$ws.Range("K141").Value = 3500
$ws.Range("L141").Value = 3500
$ws.Range("M141").Value = 3500
$ws.Range("N141").Value = "$/paquete 36 unidades"
$ws.Range("O141").Value = "Región Metropolitana"
$ws.Range("P141").Value = 97
$ws.Range("Q141").Value = 36
$ws.Range("R141").Value = "Hortaliza"
